$wb = $excel.ActiveWorkbook

# ---- helper: apply the blue-underline "HyperLink" look used elsewhere in
# this workbook (style index 1: font with Underline + color FF6495ED) ----
function Set-HyperlinkLook($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = 15570276   # RGB(0x64,0x95,0xED) == FF6495ED (BGR-packed OLE color)
}

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the "Ready for handoff" status text becomes
# "Handed back: in sync with en-US" for both tracked files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet: status text, new "Latest Target File" / "Latest Handback
# File" hyperlink cells, and the "Latest Handback DateTime" stamps.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsZh.Range("E2").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce4ae053c47528705576bad980eafdba0dfc3402/e2e/0686809e-566d-475d-8c26-c940e58fc9a0.md", "", "", "0686809e-566d-475d-8c26-c940e58fc9a0.md") | Out-Null
Set-HyperlinkLook $wsZh.Range("E2")

$wsZh.Range("F2").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48257adda512378ed4aeb0fdd3b38277936f3d21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf", "", "", "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf") | Out-Null
Set-HyperlinkLook $wsZh.Range("F2")

$wsZh.Range("E3").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce4ae053c47528705576bad980eafdba0dfc3402/e2e/f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md", "", "", "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md") | Out-Null
Set-HyperlinkLook $wsZh.Range("E3")

$wsZh.Range("F3").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/48257adda512378ed4aeb0fdd3b38277936f3d21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf", "", "", "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf") | Out-Null
Set-HyperlinkLook $wsZh.Range("F3")

$wsZh.Range("G2").Value = "2016-03-10 14:39:16"
$wsZh.Range("G3").Value = "2016-03-10 14:39:16"

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, different datetime stamp.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$wsDe.Range("E2").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ce4ae053c47528705576bad980eafdba0dfc3402/e2e/0686809e-566d-475d-8c26-c940e58fc9a0.md", "", "", "0686809e-566d-475d-8c26-c940e58fc9a0.md") | Out-Null
Set-HyperlinkLook $wsDe.Range("E2")

$wsDe.Range("F2").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b38b617ce508a5d3521813e60a6a345e857880f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf", "", "", "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf") | Out-Null
Set-HyperlinkLook $wsDe.Range("F2")

$wsDe.Range("E3").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ce4ae053c47528705576bad980eafdba0dfc3402/e2e/f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md", "", "", "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md") | Out-Null
Set-HyperlinkLook $wsDe.Range("E3")

$wsDe.Range("F3").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b38b617ce508a5d3521813e60a6a345e857880f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf", "", "", "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf") | Out-Null
Set-HyperlinkLook $wsDe.Range("F3")

$wsDe.Range("G2").Value = "2016-03-10 14:39:23"
$wsDe.Range("G3").Value = "2016-03-10 14:39:23"
